# Updates cryptocurrency price/volume figures in the "cryptos" worksheet.
# Each entry below corresponds to a refreshed quote pulled from the source
# feed: column D holds the latest Price text and column E the Volume(1h)
# percentage text. Values are written with a leading single-quote so Excel
# stores them as literal text (matching the original inline-string layout)
# instead of re-interpreting punctuation such as "30.710.69" as a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Cell="D2"; Value='30.710.69'},
    @{Cell="E2"; Value='  +1.73%  '},
    @{Cell="D3"; Value='1.899.71'},
    @{Cell="E3"; Value='  +2.64%  '},
    @{Cell="D4"; Value='0.9999'},
    @{Cell="E4"; Value='  +0.02%  '},
    @{Cell="D5"; Value='238.94'},
    @{Cell="E5"; Value='  +1.09%  '},
    @{Cell="D6"; Value='0.9998'},
    @{Cell="E6"; Value='  +0.00%  '},
    @{Cell="D7"; Value='0.4809'},
    @{Cell="E7"; Value='  +0.89%  '},
    @{Cell="D8"; Value='0.2847'},
    @{Cell="E8"; Value='  +0.95%  '},
    @{Cell="D9"; Value='0.06562'},
    @{Cell="E9"; Value='  +1.36%  '},
    @{Cell="D10"; Value='1.912.04'},
    @{Cell="E10"; Value='  +3.23%  '},
    @{Cell="D11"; Value='0.07460'},
    @{Cell="E11"; Value='  +2.18%  '},
    @{Cell="D12"; Value='16.75'},
    @{Cell="E12"; Value='  +2.47%  '},
    @{Cell="D13"; Value='5.121'},
    @{Cell="E13"; Value='  -0.24%  '},
    @{Cell="D14"; Value='88.22'},
    @{Cell="E14"; Value='  +1.16%  '},
    @{Cell="D15"; Value='0.6678'},
    @{Cell="E15"; Value='  +3.54%  '},
    @{Cell="D16"; Value='30.695.24'},
    @{Cell="E16"; Value='  +1.91%  '},
    @{Cell="D17"; Value='13.33'},
    @{Cell="E17"; Value='  +0.70%  '},
    @{Cell="D18"; Value='0.9996'},
    @{Cell="D19"; Value='0.000007624'},
    @{Cell="E19"; Value='  -0.02%  '},
    @{Cell="D20"; Value='232.27'},
    @{Cell="E20"; Value='  +6.44%  '},
    @{Cell="D21"; Value='2.165.02'},
    @{Cell="E21"; Value='  +3.13%  '},
    @{Cell="D22"; Value='5.305'},
    @{Cell="E22"; Value='  +0.64%  '},
    @{Cell="D23"; Value='1.000'},
    @{Cell="E23"; Value='  +0.05%  '},
    @{Cell="D24"; Value='6.237'},
    @{Cell="D25"; Value='170.03'},
    @{Cell="E25"; Value='  +3.80%  '},
    @{Cell="E26"; Value='  +1.46%  '},
    @{Cell="D27"; Value='18.73'},
    @{Cell="E27"; Value='  +1.97%  '},
    @{Cell="D28"; Value='1.971'},
    @{Cell="E28"; Value='  +3.04%  '},
    @{Cell="D29"; Value='1.401'},
    @{Cell="E29"; Value='  -1.95%  '},
    @{Cell="D30"; Value='0.1007'},
    @{Cell="E30"; Value='  +9.51%  '},
    @{Cell="D31"; Value='4.364'},
    @{Cell="E31"; Value='  +2.87%  '},
    @{Cell="D32"; Value='4.040'},
    @{Cell="E32"; Value='  +1.94%  '},
    @{Cell="D33"; Value='0.05101'},
    @{Cell="E33"; Value='  +1.67%  '},
    @{Cell="E34"; Value='  +7.14%  '},
    @{Cell="D35"; Value='0.7580'},
    @{Cell="E35"; Value='  +2.28%  '},
    @{Cell="D36"; Value='2.712'},
    @{Cell="E36"; Value='  +0.99%  '},
    @{Cell="D37"; Value='0.01880'},
    @{Cell="E37"; Value='  +2.59%  '},
    @{Cell="D38"; Value='2.661'},
    @{Cell="E38"; Value='  +1.96%  '},
    @{Cell="D39"; Value='0.9210'},
    @{Cell="E39"; Value='  +2.10%  '},
    @{Cell="D40"; Value='2.085'},
    @{Cell="E40"; Value='  +1.18%  '},
    @{Cell="D41"; Value='107.07'},
    @{Cell="E41"; Value='  +0.58%  '},
    @{Cell="D42"; Value='0.4312'},
    @{Cell="E42"; Value='  +1.52%  '},
    @{Cell="E43"; Value='  +0.63%  '},
    @{Cell="D44"; Value='5.766'},
    @{Cell="E44"; Value='  -2.73%  '},
    @{Cell="D45"; Value='7.449'},
    @{Cell="E45"; Value='  +0.29%  '},
    @{Cell="D46"; Value='64.41'},
    @{Cell="E46"; Value='  +0.98%  '},
    @{Cell="D47"; Value='0.1275'},
    @{Cell="E47"; Value='  -2.54%  '},
    @{Cell="D48"; Value='1.488'},
    @{Cell="E48"; Value='  -4.19%  '},
    @{Cell="D49"; Value='9.032'},
    @{Cell="E49"; Value='  +2.93%  '},
    @{Cell="D50"; Value='33.91'},
    @{Cell="E50"; Value='  -1.00%  '},
    @{Cell="D51"; Value='0.05675'},
    @{Cell="E51"; Value='  +0.01%  '}
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = "'" + $u.Value
}
